$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update MACRO_SCORE (column N) values for rows 2-6 from 51.15965480231979 to 51.05762969290213
$ws.Range("N2:N6").Value = 51.05762969290213
